$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at E:F, shifting old E:O (Venue..location_geom) to G:Q
$ws.Columns("E:F").Insert()

# Header row: rename/retitle the three leading index columns
$ws.Range("D1").Value = "Unnamed: 0.2"
$ws.Range("E1").Value = "Unnamed: 0.1"
$ws.Range("F1").Value = "Unnamed: 0"

# Data rows 2-20: populate the two new index columns with the same value as column D
for ($r = 2; $r -le 20; $r++) {
    $idx = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 5).Value = $idx
    $ws.Cells.Item($r, 6).Value = $idx
}

# Row 21 ("Ward"): rename the venue and clear out all of the per-venue detail columns (D:Q)
$ws.Range("B21").Value = "Ward Community Center"
$ws.Range("D21:Q21").ClearContents()
